# Re-order the three timeline rows covering 1838-1845:
#   - Row 5 (Texas annexation, previously dated/placed at 1836) moves down
#     to row 7 and its year is corrected to 1845.
#   - Row 6 (LSR's birth) moves up to row 5 (keeps its C/D/E layout).
#   - Row 7 (Shapely Ross fleeing to Texas) moves up to row 6 (keeps its
#     C/D/E layout).
# Row heights follow the relocated content, and the sheet selection is
# reset to reflect the new editing position (B6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the cell contents -------------------------------------------------
# Stage the "Texas is annexed" rich-text cell (B5) out of the way first so it
# isn't clobbered while rows 6/7 shift upward. A far-away scratch cell is
# used as a holding spot, then cleared once its contents are moved home.
$ws.Range("B5").Cut($ws.Range("B1000")) | Out-Null
$ws.Range("B5").Clear() | Out-Null

# Rows 6 and 7 (LSR born / Shapely Ross runs to Texas) each shift up one row.
$ws.Range("C6:E6").Cut($ws.Range("C5")) | Out-Null
$ws.Range("C6:E6").Clear() | Out-Null

$ws.Range("C7:E7").Cut($ws.Range("C6")) | Out-Null
$ws.Range("C7:E7").Clear() | Out-Null

# Bring the annexation text back down into its new home, row 7.
$ws.Range("B1000").Cut($ws.Range("B7")) | Out-Null
$ws.Range("B1000").Clear() | Out-Null

# --- Fix up the year column --------------------------------------------------
$ws.Range("A5").Value = 1838
$ws.Range("A6").Value = 1838
$ws.Range("A7").Value = 1845

# --- Restore wrap formatting lost by the cut/paste on the relocated cells --
$ws.Range("C5:E6").WrapText = $true
$ws.Range("B7").WrapText = $true

# --- Row heights (match the relocated content's natural sizing) ------------
$ws.Rows(5).RowHeight = 50.4
$ws.Rows(6).RowHeight = 57.45
$ws.Rows(7).RowHeight = 35.25

# --- Selection / view state --------------------------------------------------
$ws.Range("B6").Select() | Out-Null
try {
    $excel.ActiveWindow.ScrollRow = 1
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
